# Insert a new daily price observation for Berenjena (Femacal de La Calera)
# as row 241, pushing all subsequent rows down by one (334 -> 335 rows of
# data total). The new record:
#   Fecha (D) = 44784  (2022-08-11)
#   Volumen (J) = 105
#   Precio minimo (K) = 8500
#   Precio maximo (L) = 9000
#   Precio promedio ponderado (M) = 8738
#   Origen (O) = Región de Arica y Parinacota
#   Precio $/Kg (P) = 146
# All other columns repeat the constant values shared by every row in this
# sheet (Mercado ID, Mercado, Región, Codreg, Categoría ID/Nombre, Variedad,
# Calidad, Unidad de comercialización, Kg o Unidades, Clasificación).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 241:334 down to 242:335, leaving a blank row 241 that inherits
# the formatting (incl. the date number-format on column D) of the row it
# displaces, just like Excel's native "Insert Row" command.
$ws.Rows("241:241").Insert()

# Populate the newly inserted row 241 with the new data point.
$ws.Cells.Item(241, 1).Value = 3
$ws.Cells.Item(241, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(241, 3).Value = "Coquimbo"
$ws.Cells.Item(241, 4).Value = 44784
$ws.Cells.Item(241, 5).Value = 5
$ws.Cells.Item(241, 6).Value = 100112001
$ws.Cells.Item(241, 7).Value = "Berenjena"
$ws.Cells.Item(241, 8).Value = "Sin especificar"
$ws.Cells.Item(241, 9).Value = "Primera"
$ws.Cells.Item(241, 10).Value = 105
$ws.Cells.Item(241, 11).Value = 8500
$ws.Cells.Item(241, 12).Value = 9000
$ws.Cells.Item(241, 13).Value = 8738
$ws.Cells.Item(241, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(241, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(241, 16).Value = 146
$ws.Cells.Item(241, 17).Value = 60
$ws.Cells.Item(241, 18).Value = "Hortaliza"
